# Updates cryptos list values (price/volume) to match the latest scrape.
# Two coin pairs also swapped rows (15<->16, 44<->45) in the source ranking.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.344.07"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "1.818.81"
$ws.Range("E3").Value = "  -0.56%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "'315.25"
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("D6").Value = "'1.004"
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("D7").Value = "'0.5246"
$ws.Range("E7").Value = "  +2.27%  "
$ws.Range("D8").Value = "'0.3850"
$ws.Range("E8").Value = "  -2.06%  "
$ws.Range("D9").Value = "'0.08086"
$ws.Range("E9").Value = "  +5.27%  "
$ws.Range("D10").Value = "'41.84"
$ws.Range("E10").Value = "  +0.07%  "
$ws.Range("D11").Value = "'1.114"
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("D12").Value = "'6.404"
$ws.Range("D13").Value = "'1.004"
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("D14").Value = "'20.89"
$ws.Range("E14").Value = "  -1.00%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'7.412"
$ws.Range("E15").Value = "  -1.71%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "1.822.04"
$ws.Range("E16").Value = "  -0.33%  "
$ws.Range("D17").Value = "'93.92"
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("D18").Value = "'0.00001102"
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("D19").Value = "'0.06630"
$ws.Range("E19").Value = "  -1.35%  "
$ws.Range("E20").Value = "  -0.25%  "
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("D22").Value = "'6.017"
$ws.Range("E22").Value = "  -2.35%  "
$ws.Range("D23").Value = "28.385.11"
$ws.Range("E23").Value = "  -0.40%  "
$ws.Range("D24").Value = "'11.37"
$ws.Range("E24").Value = "  +1.87%  "
$ws.Range("D25").Value = "'2.243"
$ws.Range("E25").Value = "  -0.70%  "
$ws.Range("D26").Value = "'159.33"
$ws.Range("E26").Value = "  +1.47%  "
$ws.Range("D27").Value = "'20.82"
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").Value = "2.026.55"
$ws.Range("E28").Value = "  -0.60%  "
$ws.Range("D29").Value = "'2.398"
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("D30").Value = "'124.41"
$ws.Range("E30").Value = "  -0.30%  "
$ws.Range("D31").Value = "'0.1103"
$ws.Range("E31").Value = "  +1.46%  "
$ws.Range("D32").Value = "'1.079"
$ws.Range("E32").Value = "  -3.55%  "
$ws.Range("D33").Value = "'5.672"
$ws.Range("E33").Value = "  -0.12%  "
$ws.Range("D34").Value = "'3.676"
$ws.Range("E34").Value = "  +0.46%  "
$ws.Range("D35").Value = "'0.07388"
$ws.Range("E35").Value = "  +4.91%  "
$ws.Range("D36").Value = "'12.30"
$ws.Range("E36").Value = "  +9.59%  "
$ws.Range("D37").Value = "'0.2196"
$ws.Range("E37").Value = "  -0.92%  "
$ws.Range("D38").Value = "'0.02340"
$ws.Range("E38").Value = "  +0.63%  "
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("D40").Value = "'8.721"
$ws.Range("E40").Value = "  -2.88%  "
$ws.Range("D41").Value = "'0.6326"
$ws.Range("E41").Value = "  +0.74%  "
$ws.Range("D42").Value = "'1.182"
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D43").Value = "'1.380"
$ws.Range("E43").Value = "  -0.86%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'13.44"
$ws.Range("E44").Value = "  -0.56%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.6125"
$ws.Range("E45").Value = "  +3.73%  "
$ws.Range("D46").Value = "'3.784"
$ws.Range("E46").Value = "  +1.75%  "
$ws.Range("D47").Value = "'127.21"
$ws.Range("E47").Value = "  +1.79%  "
$ws.Range("D48").Value = "'1.981"
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("D50").Value = "'0.06885"
$ws.Range("E50").Value = "  -0.64%  "
$ws.Range("D51").Value = "'1.066"
$ws.Range("E51").Value = "  -0.31%  "
